$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.456404
$ws.Range("H2").Value = 1.369212
$ws.Range("I2").Value = 0.01914960767004715
$ws.Range("J2").Value = 0.01914960767004715
$ws.Range("M2").Value = 2.294987
$ws.Range("N2").Value = 6.884961000000001
$ws.Range("O2").Value = 0.0158275801650097
$ws.Range("P2").Value = 0.0158275801650097
$ws.Range("Q2").Value = 1.047441246748
$ws.Range("R2").Value = 9.426971220732002
$ws.Range("S2").Value = 0.000303091950526156
$ws.Range("T2").Value = 0.000303091950526156

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.456404
$ws.Range("H3").Value = 1.369212
$ws.Range("I3").Value = 0.01914960767004715
$ws.Range("J3").Value = 0.01914960767004715
$ws.Range("O3").Value = 0.769602070219672
$ws.Range("P3").Value = 0.7696020702196722
$ws.Range("Q3").Value = 50.93090311510933
$ws.Range("R3").Value = 458.378128035984
$ws.Range("S3").Value = 0.0147375777067628
$ws.Range("T3").Value = 0.0147375777067628

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.456404
$ws.Range("H4").Value = 1.369212
$ws.Range("I4").Value = 0.01914960767004715
$ws.Range("J4").Value = 0.01914960767004715
$ws.Range("M4").Value = 31.11253633333333
$ws.Range("N4").Value = 93.337609
$ws.Range("O4").Value = 0.2145703496153182
$ws.Range("P4").Value = 0.2145703496153182
$ws.Range("Q4").Value = 14.19988603267867
$ws.Range("R4").Value = 127.798974294108
$ws.Range("S4").Value = 0.004108938012758197
$ws.Range("T4").Value = 0.004108938012758197

# Row 5
$ws.Range("I5").Value = 0.8285024587002443
$ws.Range("J5").Value = 0.8285024587002443
$ws.Range("M5").Value = 2.294987
$ws.Range("N5").Value = 6.884961000000001
$ws.Range("O5").Value = 0.0158275801650097
$ws.Range("P5").Value = 0.0158275801650097
$ws.Range("Q5").Value = 45.31725470449968
$ws.Range("R5").Value = 407.8552923404971
$ws.Range("S5").Value = 0.01311318908198576
$ws.Range("T5").Value = 0.01311318908198576

# Row 6
$ws.Range("I6").Value = 0.8285024587002443
$ws.Range("J6").Value = 0.8285024587002443
$ws.Range("O6").Value = 0.769602070219672
$ws.Range("P6").Value = 0.7696020702196722
$ws.Range("S6").Value = 0.6376172073977963
$ws.Range("T6").Value = 0.6376172073977965

# Row 7
$ws.Range("I7").Value = 0.8285024587002443
$ws.Range("J7").Value = 0.8285024587002443
$ws.Range("M7").Value = 31.11253633333333
$ws.Range("N7").Value = 93.337609
$ws.Range("O7").Value = 0.2145703496153182
$ws.Range("P7").Value = 0.2145703496153182
$ws.Range("S7").Value = 0.1777720622204622
$ws.Range("T7").Value = 0.1777720622204622

# Row 8
$ws.Range("I8").Value = 0.1523479336297086
$ws.Range("J8").Value = 0.1523479336297086
$ws.Range("M8").Value = 2.294987
$ws.Range("N8").Value = 6.884961000000001
$ws.Range("O8").Value = 0.0158275801650097
$ws.Range("P8").Value = 0.0158275801650097
$ws.Range("Q8").Value = 8.333095502013002
$ws.Range("R8").Value = 74.99785951811701
$ws.Range("S8").Value = 0.002411299132497791
$ws.Range("T8").Value = 0.002411299132497791

# Row 9
$ws.Range("I9").Value = 0.1523479336297086
$ws.Range("J9").Value = 0.1523479336297086
$ws.Range("O9").Value = 0.769602070219672
$ws.Range("P9").Value = 0.7696020702196722
$ws.Range("S9").Value = 0.117247285115113
$ws.Range("T9").Value = 0.117247285115113

# Row 10
$ws.Range("I10").Value = 0.1523479336297086
$ws.Range("J10").Value = 0.1523479336297086
$ws.Range("M10").Value = 31.11253633333333
$ws.Range("N10").Value = 93.337609
$ws.Range("O10").Value = 0.2145703496153182
$ws.Range("P10").Value = 0.2145703496153182
$ws.Range("S10").Value = 0.03268934938209787
$ws.Range("T10").Value = 0.03268934938209788
